$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.551.35"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "1.659.79"
$ws.Range("E3").Value = "  -3.37%  "
$ws.Range("E4").Value = "  +0.99%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("E6").Value = "  -1.29%  "
$ws.Range("E7").Value = "  +0.98%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.27"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.73%  "
$ws.Range("E9").Value = "  -2.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0619"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("E11").Value = "  -2.10%  "
$ws.Range("D12").Value = "1.892.98"
$ws.Range("E12").Value = "  -3.18%  "
$ws.Range("D13").Value = "1.659.56"
$ws.Range("E13").Value = "  -3.04%  "
$ws.Range("E14").Value = "  -3.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.549"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.35%  "
$ws.Range("E16").Value = "  -3.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "245.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("D18").Value = "27.536.00"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").Value = "0.0₃0728"
$ws.Range("E19").Value = "  -3.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.18%  "
$ws.Range("E21").Value = "  +0.90%  "
$ws.Range("E22").Value = "  -3.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.16%  "
$ws.Range("E24").Value = "  -3.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.71%  "
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.112"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.94%  "
$ws.Range("E30").Value = "  +4.95%  "
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("E32").Value = "  -2.39%  "
$ws.Range("D33").Value = "1.436.03"
$ws.Range("E33").Value = "  -7.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.03%  "
$ws.Range("E35").Value = "  -8.40%  "
$ws.Range("E36").Value = "  -0.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.931"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.35%  "
$ws.Range("E38").Value = "  -5.83%  "
$ws.Range("E39").Value = "  -2.67%  "
$ws.Range("E40").Value = "  -1.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.22%  "
$ws.Range("E42").Value = "  +0.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.790"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "1.801.80"
$ws.Range("E45").Value = "  -2.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.59%  "
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.61%  "
$ws.Range("E49").Value = "  +4.82%  "
$ws.Range("E50").Value = "  -4.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.50%  "
